$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Onion Red -> Current Quantity 0 -> 1, Requested quantity 4 -> 0
# Column B stores its quantities as text. Copy a cell that already holds the
# text value "1" (B2) into B3 via Paste Special Values so the new value is
# written as text (matching the sheet's existing text-based column B),
# without introducing any new cell style.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C3").Value = 0

# Row 68: Thai Chilli -> Current Quantity cleared (was 2), Requested quantity 6 -> 0
$ws.Range("B68").ClearContents()
$ws.Range("C68").Value = 0

# Row 69: Sattu -> Current Quantity cleared (was 0), Requested quantity 3 -> 0
$ws.Range("B69").ClearContents()
$ws.Range("C69").Value = 0

# Row 134: Very Spicy Red Chilli Powder -> Current Quantity cleared (was 1), Requested quantity 4 -> 0
$ws.Range("B134").ClearContents()
$ws.Range("C134").Value = 0
